$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Fill in the new agenda entries (rows 8-11)
# Columns: A=Técnico, B=ID, C=Cliente, D=Descricao, G=Status

$ws.Range("A8").Value = "Roberto"
$ws.Range("B8").Value = "'0869"
$ws.Range("C8").Value = "Recapel"
$ws.Range("D8").Value = "Zona aberta, preventiva essa o.s."
$ws.Range("G8").Value = "Pendente"

$ws.Range("A9").Value = "Roberto"
$ws.Range("B9").Value = "'0706"
$ws.Range("C9").Value = "Lar das Meninas"
$ws.Range("D9").Value = "Várias câmeras sem imagem."
$ws.Range("G9").Value = "Pendente"

$ws.Range("A10").Value = "Roberto"
$ws.Range("B10").Value = "'0773"
$ws.Range("C10").Value = "Escola Antônio Gonçalves de Matos"
$ws.Range("D10").Value = "Várias câmeras sem imagem e acesso remoto pra Cida."
$ws.Range("G10").Value = "Pendente"

$ws.Range("A11").Value = "Roberto"
$ws.Range("B11").Value = "'0304"
$ws.Range("C11").Value = "Cimentão"
$ws.Range("D11").Value = "Ordem consta disparos frequentes em algumas zonas."
$ws.Range("G11").Value = "Pendente"

# Update view state: scroll position and active cell selection
$excel.Goto($ws.Range("F1"), $true)
$ws.Range("H11").Select()
